$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be misread as a number by Excel's
# COM Value coercion (e.g. "1.00" -> 1, "0.520" -> 0.52). Force these
# specific cells to Text format first so the literal string round-trips.
$textCells = @("D5","D6","D7","D11","D12","D13","D14","D17","D19","D21","D22","D23","D24","D25","D27","D28","D30","D31","D32","D35","D36","D37","D38","D41","D42","D44","D45","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "66.652.88"
$ws.Range("E2").Value = "  +2.63%  "

# Row 3
$ws.Range("D3").Value = "3.203.51"
$ws.Range("E3").Value = "  +1.45%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "602.14"
$ws.Range("E5").Value = "  +3.73%  "

# Row 6
$ws.Range("D6").Value = "156.10"
$ws.Range("E6").Value = "  +4.38%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  +6.39%  "

# Row 9
$ws.Range("D9").Value = "3.201.11"
$ws.Range("E9").Value = "  +1.39%  "

# Row 10
$ws.Range("E10").Value = "  +1.58%  "

# Row 11
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  -4.68%  "

# Row 12
$ws.Range("D12").Value = "0.520"
$ws.Range("E12").Value = "  +3.80%  "

# Row 13
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  +2.05%  "

# Row 14
$ws.Range("D14").Value = "39.27"
$ws.Range("E14").Value = "  +5.57%  "

# Row 15
$ws.Range("D15").Value = "3.729.11"
$ws.Range("E15").Value = "  +1.43%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.673.37"
$ws.Range("E16").Value = "  +2.68%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "7.51"
$ws.Range("E17").Value = "  +5.04%  "

# Row 18
$ws.Range("D18").Value = "3.201.97"
$ws.Range("E18").Value = "  +1.57%  "

# Row 19
$ws.Range("D19").Value = "526.58"
$ws.Range("E19").Value = "  +4.37%  "

# Row 20
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("D21").Value = "15.51"
$ws.Range("E21").Value = "  +3.94%  "

# Row 22
$ws.Range("D22").Value = "0.742"
$ws.Range("E22").Value = "  +3.73%  "

# Row 23
$ws.Range("D23").Value = "8.17"
$ws.Range("E23").Value = "  +5.50%  "

# Row 24
$ws.Range("D24").Value = "15.02"
$ws.Range("E24").Value = "  -1.05%  "

# Row 25
$ws.Range("D25").Value = "85.87"
$ws.Range("E25").Value = "  +1.55%  "

# Row 26
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +2.86%  "

# Row 28
$ws.Range("D28").Value = "3.02"
$ws.Range("E28").Value = "  +2.97%  "

# Row 29
$ws.Range("E29").Value = "  +8.87%  "

# Row 30
$ws.Range("D30").Value = "2.95"
$ws.Range("E30").Value = "  +5.67%  "

# Row 31
$ws.Range("D31").Value = "7.04"
$ws.Range("E31").Value = "  +10.54%  "

# Row 32
$ws.Range("D32").Value = "28.37"
$ws.Range("E32").Value = "  +2.72%  "

# Row 33
$ws.Range("E33").Value = "  +3.13%  "

# Row 34
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("D35").Value = "6.57"
$ws.Range("E35").Value = "  +1.24%  "

# Row 36
$ws.Range("D36").Value = "510.78"
$ws.Range("E36").Value = "  +6.19%  "

# Row 37
$ws.Range("D37").Value = "54.89"
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").Value = "0.0908"
$ws.Range("E38").Value = "  +1.58%  "

# Row 40
$ws.Range("E40").Value = "  +8.86%  "

# Row 41
$ws.Range("D41").Value = "8.93"
$ws.Range("E41").Value = "  +2.08%  "

# Row 42
$ws.Range("D42").Value = "2.89"
$ws.Range("E42").Value = "  -1.13%  "

# Row 43
$ws.Range("E43").Value = "  +15.48%  "

# Row 44
$ws.Range("D44").Value = "0.301"
$ws.Range("E44").Value = "  +6.35%  "

# Row 45
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").Value = "  +1.04%  "

# Row 46
$ws.Range("D46").Value = "2.898.77"
$ws.Range("E46").Value = "  -3.15%  "

# Row 47
$ws.Range("E47").Value = "  +1.11%  "

# Row 48
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  +9.86%  "

# Row 49
$ws.Range("E49").Value = "  +3.75%  "

# Row 50
$ws.Range("D50").Value = "2.36"
$ws.Range("E50").Value = "  +5.25%  "

# Row 51
$ws.Range("E51").Value = "  -0.02%  "

